$d = $word.ActiveDocument

# 1. Update the letter date (September 19 -> September 21).
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 1) | Out-Null

# 2. Split the single "street, city state zip" mailing-address line into two
#    separate paragraphs: "2654 Greenrock Road" and "Milpitas, CA 95035".
#    Only the first occurrence (the recipient's mailing address block) is
#    affected; the "PROPERTY ADDRESS:" line further down keeps the
#    combined text.
$rng = $d.Content
$rng.Find.Execute("2654 Greenrock Road, Milpitas CA 95035", $true, $false, $false, $false, $false,
                   $true, 1, $false, "2654 Greenrock Road", 1) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.Move(4, 1) | Out-Null
$rng.Text = "Milpitas, CA 95035"

# 3. Remove the now-redundant blank "NoSpacing" paragraph that used to sit
#    immediately after the "...Board of Directors" signature line.
$found = $d.Content
$found.Find.Execute("Board of Directors") | Out-Null
$para = $found.Paragraphs(1)
$next = $para.Next()
if ($next -ne $null -and $next.Range.Text -eq "`r") {
    $next.Range.Delete() | Out-Null
}
